$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update the label text in column A (rows 29-36) ---
# The underlying shared-string table is being cleaned up / renamed:
#   NFTINDCORPGRPITG25%C      -> NFTINDCORPGRPITG25PC   (moved to row 29)
#   NFT500MULCINDMFG50:30:20  -> NFT500MULCINDMFG50_30_20 (moved to row 34)
#   NFT500MULCINFS50:30:20    -> NFT500MULCINFS50_30_20   (moved to row 35)
#   NFTSMEEMG moves up to row 36
# and the remaining labels keep their text, just shift rows.
$ws.Range("A29").Value = "NFTINDCORPGRPITG25PC"
$ws.Range("A30").Value = "NFTTRANSLOG"
$ws.Range("A31").Value = "NFT100L15"
$ws.Range("A32").Value = "NFT50SH"
$ws.Range("A33").Value = "NFT500SH"
$ws.Range("A34").Value = "NFT500MULCINDMFG50_30_20"
$ws.Range("A35").Value = "NFT500MULCINFS50_30_20"
$ws.Range("A36").Value = "NFTSMEEMG"

# --- Widen column A to fit the new (longer) labels ---
$ws.Columns.Item(1).ColumnWidth = 25.5

# --- Update the view: scroll so row 11 is at the top, and select E29 ---
$ws.Range("A11").Select()
$excel.ActiveWindow.ScrollRow = 11
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("E29").Select()
